# "fixing key problems in old database sheets"
#
# - bioSampleNumber (column C) for data rows 2-45 gets shifted by +206
#   (old values 1-44 become 207-250).
# - Data rows 2-45 get a tighter row height (13.8pt instead of the
#   default 15pt).
# - The view is scrolled/selected to show the refreshed bioSampleNumber
#   column: C2:C45 selected with C2 as the active cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bump bioSampleNumber values in column C for rows 2..45 by 206.
for ($r = 2; $r -le 45; $r++) {
    $ws.Cells.Item($r, 3).Value = $r + 205
}

# Tighten row height for the data rows (2..45) to 13.8pt.
$ws.Range("A2:A45").RowHeight = 13.8

# Reflect the edited column in the sheet's selection.
[void]$ws.Range("C2:C45").Select()
